$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new "Result" value for the first data row (N2) as described
# in the commit: a new "Fail" result is recorded for the first test case.
$ws.Range("N2").Value = "Fail"
